# Apply "update: att product started" change:
#  - Row 13 gets real product data (A13/B13 newly filled, C13/D13 text
#    changed, E13/F13/H13 become text-typed numeric strings, G13 becomes 22)
#  - Column I (the stray "Unnamed: 8" column) is removed entirely

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13 updates -------------------------------------------------------

$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "1234"
$ws.Range("A13").Style = "Normal"

$ws.Range("B13").Value = "OI"

$ws.Range("C13").Value = "Oi"

$ws.Range("D13").Value = "Oi"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "112"
$ws.Range("E13").Style = "Normal"

$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "134"
$ws.Range("F13").Style = "Normal"

$ws.Range("G13").Value = 22

$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "4"
$ws.Range("H13").Style = "Normal"

# --- Drop the stray column I ("Unnamed: 8") --------------------------------

$ws.Columns("I").Delete()
